# Script to use same logic as bulk_insert.py, but get fragments in reactants and products
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Append brand-new fragment rows at the bottom of the sheet first
#    (Acetylenic Carbon, Anhydride, Amidinium, Cyanamide) - this is the order
#    in which the corresponding shared strings were originally created.
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Acetylenic Carbon"
$ws.Range("B23").Value = "[`$([CX2]#C)]"

$ws.Range("A24").Value = "Anhydride"
$ws.Range("B24").Value = "[CX3](=[OX1])[OX2][CX3](=[OX1])"

$ws.Range("A25").Value = "Amidinium"
$ws.Range("B25").Value = "[NX3][CX3]=[NX3+]"

$ws.Range("A26").Value = "Cyanamide"
$ws.Range("B26").Value = "[NX3][CX2]#[NX1]"

# ---------------------------------------------------------------------------
# 2) Insert a brand-new row 7 (pushing Aldehyde and everything below it down
#    by one row), and populate it with a Sub-Fragment/Sub-SMARTS pair
#    (Enamine) followed by the main Fragment/SMARTS pair (Aniline Nitrogen).
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).Insert()

$ws.Range("C7").Value = "Enamine"
$ws.Range("D7").Value = "[NX3][CX3]=[CX3]"

$ws.Range("A7").Value = "Aniline Nitrogen"
$ws.Range("B7").Value = "[NX3][`$(C=C),`$(cc)]"

# ---------------------------------------------------------------------------
# 3) Append more new fragment rows (Azide, Azo Nitrogen, Sulfuric acid,
#    Alkyl Carbon). For the first three the SMARTS (column B) was entered
#    before the fragment name (column A).
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = "[`$(*-[NX2-]-[NX2+]#[NX1]),`$(*-[NX2]=[NX2+]=[NX1-])]"
$ws.Range("A28").Value = "Azide"

$ws.Range("B29").Value = "[NX2]=N"
$ws.Range("A29").Value = "Azo Nitrogen"

$ws.Range("B30").Value = "[`$([SX4](=O)(=O)(O)O),`$([SX4+2]([O-])([O-])(O)O)]"
$ws.Range("A30").Value = "Sulfuric acid"

$ws.Range("A31").Value = "Alkyl Carbon"
$ws.Range("B31").Value = "[CX4]"

# ---------------------------------------------------------------------------
# 4) Re-append a few already-existing fragments at the bottom again
#    (duplicates of Allenic Carbon, Vinylic Carbon and Acetylenic Carbon).
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = "Allenic Carbon"
$ws.Range("B32").Value = "[`$([CX2](=C)=C)]"

$ws.Range("A33").Value = "Vinylic Carbon"
$ws.Range("B33").Value = "[`$([CX3]=[CX3])]"

$ws.Range("A34").Value = "Acetylenic Carbon"
$ws.Range("B34").Value = "[`$([CX2]#C)]"

# ---------------------------------------------------------------------------
# 5) Cosmetic sheet-level changes: column widths, zoom, selection, and
#    page orientation.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.346354166666666
$ws.Columns.Item(2).ColumnWidth = 51.166666666666664
$ws.Columns.Item(3).ColumnWidth = 13.983072916666666

$ws.Application.ActiveWindow.Zoom = 106
$ws.Range("D17").Select()

$ws.PageSetup.Orientation = 1
